$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived NATMI values for rows 2-7 (columns E:T)
$values = @{
    2 = @{ E=3; F=1; G=11.629057; H=34.887171; I=0.3062678464977661; J=0.3062678464977662;
           K=3; L=1; M=1.357948666666666; N=4.073846; O=0.936754494361095; P=0.936754494361095;
           Q=15.79166244774066; R=142.124962029666; S=0.2868977816850763; T=0.2868977816850764 }
    3 = @{ E=3; F=1; G=11.629057; H=34.887171; I=0.3062678464977661; J=0.3062678464977662;
           K=2; L=0.6666666666666666; M=0.09168266666666665; N=0.275048; O=0.06324550563890498; P=0.06324550563890496;
           Q=1.066182956578666; R=9.595646609207996; S=0.01937006481268975; T=0.01937006481268975 }
    4 = @{ I=0.6269156120645606; J=0.6269156120645607;
           K=3; L=1; M=1.357948666666666; N=4.073846; O=0.936754494361095; P=0.936754494361095;
           Q=32.32477663636978; R=290.922989727328; S=0.5872660171866139; T=0.587266017186614 }
    5 = @{ I=0.6269156120645606; J=0.6269156120645607;
           K=2; L=0.6666666666666666; M=0.09168266666666665; N=0.275048; O=0.06324550563890498; P=0.06324550563890496;
           Q=2.182425443740444; R=19.641828993664; S=0.03964959487794673; T=0.03964959487794673 }
    6 = @{ G=2.537038666666667; H=7.611116; I=0.06681654143767324; J=0.06681654143767324;
           K=3; L=1; M=1.357948666666666; N=4.073846; O=0.936754494361095; P=0.936754494361095;
           Q=3.445168274681777; R=31.006514472136; S=0.06259069548940475; T=0.06259069548940475 }
    7 = @{ G=2.537038666666667; H=7.611116; I=0.06681654143767324; J=0.06681654143767324;
           K=2; L=0.6666666666666666; M=0.09168266666666665; N=0.275048; O=0.06324550563890498; P=0.06324550563890496;
           Q=0.2326024703964444; R=2.093422233568; S=0.004225845948268491; T=0.00422584594826849 }
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$row").Value = $rowVals[$col]
    }
}
